# Generate Report for Handback
# Marks the zh-cn / de-de handback rows as "in sync" (handback complete),
# refreshes the "Latest Handback DateTime" stamps, and clears the stale
# "handback file is not latest" error message now that the report is
# regenerated.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-08-27 16:48:36"
$zhcn.Range("P2").Value = ""

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-08-27 16:48:42"
$dede.Range("P2").Value = ""

# --- Column width refresh (the longer status text / now-empty error
#     column change the auto-fit widths of the affected columns) ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
